# Insert a new data row at row 17 (pushing the existing rows 17-128 down to
# 18-129) and populate it with the new weekly observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(17).Insert()

$ws.Range("A17").Value2 = 6
$ws.Range("B17").Value2 = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C17").Value2 = "Metropolitana"
$ws.Range("D17").Value2 = 44547
$ws.Range("E17").Value2 = 13
$ws.Range("F17").Value2 = "Fruta"
$ws.Range("G17").Value2 = 100101
$ws.Range("H17").Value2 = "Berries"
$ws.Range("I17").Value2 = 100101004
$ws.Range("J17").Value2 = "Frambuesa"
$ws.Range("K17").Value2 = "Sin especificar"
$ws.Range("L17").Value2 = "Primera"
$ws.Range("M17").Value2 = 250
$ws.Range("N17").Value2 = 7000
$ws.Range("O17").Value2 = 7000
$ws.Range("P17").Value2 = 7000
$ws.Range("Q17").Value2 = "$/bandeja 2 kilos"
$ws.Range("R17").Value2 = "Provincia de Curicó"
$ws.Range("S17").Value2 = 3500
$ws.Range("T17").Value2 = 2
